$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 4 ("IMPLEMENTATION AND CODE"): the Git-hub links table.
# The last row ("PPT Link") had an empty second cell; fill it in
# with the link to the uploaded .pptx, styled the same blue
# (RGB 0,112,192 / hex 0070C0) used by the other link cells.
# -----------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$tableShape = $slide4.Shapes.Item(5)
$table = $tableShape.Table
$pptLinkCell = $table.Cell(4, 2)
$pptLinkRange = $pptLinkCell.Shape.TextFrame.TextRange
$pptLinkRange.Text = "https://github.com/aswathaanalina-19/connect-4-ai/blob/main/connect-4-with-AI(ppt)-ASWA.pptx"
$pptLinkRange.Font.Color.RGB = 12611584

# Adding that run makes the "PPT Link" row taller. Resize just that
# row (EMU 2,367,280 = 186.4pt) so the table's overall extent is
# persisted as cy=3,108,960 EMU, without disturbing the other rows'
# stored heights (370,840 EMU each / header row's 0).
$pptLinkRow = $table.Rows(4)
$pptLinkRow.Height = (2367280 / 12700.0)

# -----------------------------------------------------------------
# Slide 7 ("REFERENCES"): the last GitHub-repositories bullet had
# "performance " and "tuning ." split across two runs with
# identical formatting; normalize them into a single run.
# -----------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
$refShape = $slide7.Shapes.Item(2)
$refRange = $refShape.TextFrame.TextRange
$fullText = $refRange.Text
$mergeStart = $fullText.IndexOf("performance tuning .")
$mergeRange = $refRange.Characters($mergeStart + 1, 21)
$mergeRange.Text = "performance tuning ."
